# Project is completed for now
# Add the new "Mustafa kamal" user account record as row 13 on the
# "userAccount" sheet (A1:E12 -> A1:E13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("userAccount")

$ws.Range("A13").Value = "Mustafa kamal"
$ws.Range("B13").Value = "mustafa"

# Columns C (Password) and E (Phone) hold digit-only text in this sheet
# (e.g. "123", "01521206720") - format as text first so the numeric-looking
# values are stored as strings rather than being coerced to numbers.
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "123"

$ws.Range("D13").Value = "fdsf"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "45325"
